# Adds more texture-part rows (AN/AP/AQ columns) to Sheet1, rows 40-97,
# matching the "Added more texture parts" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, AN (part index), AP (x), AQ (y) - texture part coordinates table
$textureParts = @(
    @(40, 0, 511, 383),
    @(41, 1, 447, 448),
    @(42, 2, 381, 448),
    @(43, 3, 319, 384),
    @(44, 4, 256, 447),
    @(45, 5, 318, 510),
    @(46, 6, 383, 512),
    @(47, 7, 446, 512),
    @(48, 8, 511, 512),
    @(49, 9, 575, 448),
    @(50, 10, 576, 190),
    @(51, 11, 704, 190),
    @(52, 12, 704, 127),
    @(53, 13, 383, 127),
    @(54, 14, 383, 190),
    @(55, 15, 511, 190),
    @(56, 16, 512, 703),
    @(57, 17, 447, 639),
    @(58, 18, 383, 638),
    @(59, 19, 318, 703),
    @(60, 20, 256, 650),
    @(61, 21, 319, 577),
    @(62, 22, 383, 577),
    @(63, 23, 447, 577),
    @(64, 24, 512, 577),
    @(65, 25, 575, 639),
    @(66, 26, 574, 864),
    @(67, 27, 704, 865),
    @(68, 28, 704, 928),
    @(69, 29, 383, 928),
    @(70, 30, 383, 862),
    @(71, 31, 511, 861),
    @(72, 32, 639, 255),
    @(73, 33, 639, 448),
    @(74, 34, 639, 511),
    @(75, 35, 575, 571),
    @(76, 36, 511, 575),
    @(77, 37, 319, 576),
    @(78, 38, 255, 575),
    @(79, 39, 191, 512),
    @(80, 40, 194, 384),
    @(81, 41, 255, 322),
    @(82, 42, 382, 320),
    @(83, 43, 443, 383),
    @(84, 44, 94, 736),
    @(85, 45, 157, 736),
    @(86, 46, 383, 703),
    @(87, 47, 447, 765),
    @(88, 48, 448, 384),
    @(89, 49, 448, 256),
    @(90, 50, 383, 255),
    @(91, 51, 319, 191),
    @(92, 52, 318, 126),
    @(93, 53, 383, 65),
    @(94, 54, 700, 63),
    @(95, 55, 768, 129),
    @(96, 56, 767, 191),
    @(97, 57, 703, 256)
)

foreach ($part in $textureParts) {
    $row = $part[0]
    $ws.Cells.Item($row, 40).Value = $part[1]   # AN = part index
    $ws.Cells.Item($row, 42).Value = $part[2]   # AP = x coordinate
    $ws.Cells.Item($row, 43).Value = $part[3]   # AQ = y coordinate
}

Write-Output ("Wrote " + $textureParts.Count + " texture-part rows (AN/AP/AQ), rows 40-97")
